# Updated graphs and tables from Stata
# Refresh the correlation-matrix values in WorkingFolder/Tables/corr3mvM.xlsx
# (columns B:R, skipping I, rows 2-18) with the new Stata output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.189159659647154
$ws.Range("E2").Value = 0.0956577578795192
$ws.Range("F2").Value = 0.03298253691557487
$ws.Range("G2").Value = -0.09299971502374568
$ws.Range("H2").Value = -0.05989893757236041
$ws.Range("J2").Value = -0.05824571609997801
$ws.Range("K2").Value = -0.1494911127189555
$ws.Range("L2").Value = -0.1043516007832907
$ws.Range("M2").Value = -0.09656379635806157
$ws.Range("N2").Value = 0.04217011538748306
$ws.Range("O2").Value = -0.1144752272881613
$ws.Range("P2").Value = 0.02034454368830499
$ws.Range("Q2").Value = -0.143378829685207
$ws.Range("R2").Value = 0.03120943056409502
$ws.Range("D3").Value = 0.2125301326069027
$ws.Range("E3").Value = -0.05531880183914338
$ws.Range("F3").Value = 0.002678038600984968
$ws.Range("G3").Value = 0.102593418178592
$ws.Range("H3").Value = -0.04345523801470844
$ws.Range("J3").Value = -0.02294459231763052
$ws.Range("K3").Value = 0.158775156948585
$ws.Range("L3").Value = 0.09866967570445913
$ws.Range("M3").Value = -0.04734339998996136
$ws.Range("N3").Value = -0.09583277978816673
$ws.Range("O3").Value = 0.1405959713080313
$ws.Range("P3").Value = -0.04777815539842905
$ws.Range("Q3").Value = 0.1309388383387151
$ws.Range("R3").Value = 0.01224759455026274
$ws.Range("B4").Value = -0.189159659647154
$ws.Range("C4").Value = 0.2125301326069027
$ws.Range("E4").Value = -0.7493962068288744
$ws.Range("F4").Value = 0.7694196976433442
$ws.Range("G4").Value = 0.3179411469843243
$ws.Range("H4").Value = -0.2494502847149288
$ws.Range("J4").Value = -0.3033481172167592
$ws.Range("K4").Value = 0.4015523515806572
$ws.Range("L4").Value = -0.5727500228386581
$ws.Range("M4").Value = -0.07574937763848384
$ws.Range("N4").Value = -0.4971150376806725
$ws.Range("O4").Value = 0.1114722979239853
$ws.Range("P4").Value = -0.3819033532437513
$ws.Range("Q4").Value = 0.3454478151133055
$ws.Range("R4").Value = -0.5891277964904179
$ws.Range("B5").Value = 0.0956577578795192
$ws.Range("C5").Value = -0.05531880183914338
$ws.Range("D5").Value = -0.7493962068288744
$ws.Range("F5").Value = -0.8954460515141784
$ws.Range("G5").Value = -0.3281901268206336
$ws.Range("H5").Value = 0.3748955531446156
$ws.Range("J5").Value = 0.3816964293772955
$ws.Range("K5").Value = -0.6093742791727917
$ws.Range("L5").Value = 0.6801099074534445
$ws.Range("M5").Value = 0.07721240671536068
$ws.Range("N5").Value = 0.4942970076750814
$ws.Range("O5").Value = -0.02813934050008031
$ws.Range("P5").Value = 0.3441744846429838
$ws.Range("Q5").Value = -0.550566853378591
$ws.Range("R5").Value = 0.6426805496216214
$ws.Range("B6").Value = 0.03298253691557487
$ws.Range("C6").Value = 0.002678038600984968
$ws.Range("D6").Value = 0.7694196976433442
$ws.Range("E6").Value = -0.8954460515141784
$ws.Range("G6").Value = 0.3584452866038336
$ws.Range("H6").Value = -0.3817243282353112
$ws.Range("J6").Value = -0.4261300346409569
$ws.Range("K6").Value = 0.3936451685116808
$ws.Range("L6").Value = -0.694086968412218
$ws.Range("M6").Value = -0.09685851345657408
$ws.Range("N6").Value = -0.5526442553679262
$ws.Range("O6").Value = 0.05383151016975546
$ws.Range("P6").Value = -0.4195533350460475
$ws.Range("Q6").Value = 0.3596911983350313
$ws.Range("R6").Value = -0.5916236419409803
$ws.Range("B7").Value = -0.09299971502374568
$ws.Range("C7").Value = 0.102593418178592
$ws.Range("D7").Value = 0.3179411469843243
$ws.Range("E7").Value = -0.3281901268206336
$ws.Range("F7").Value = 0.3584452866038336
$ws.Range("H7").Value = 0.2619152394178419
$ws.Range("J7").Value = 0.2529850742589501
$ws.Range("K7").Value = 0.3350808734353892
$ws.Range("L7").Value = -0.08643860460969868
$ws.Range("M7").Value = 0.779448656064079
$ws.Range("N7").Value = 0.05749213078032143
$ws.Range("O7").Value = 0.1739025695239392
$ws.Range("P7").Value = 0.02311440040280032
$ws.Range("Q7").Value = 0.2711932158332709
$ws.Range("R7").Value = -0.1406704254457551
$ws.Range("B8").Value = -0.05989893757236041
$ws.Range("C8").Value = -0.04345523801470844
$ws.Range("D8").Value = -0.2494502847149288
$ws.Range("E8").Value = 0.3748955531446156
$ws.Range("F8").Value = -0.3817243282353112
$ws.Range("G8").Value = 0.2619152394178419
$ws.Range("J8").Value = 0.9668095370549465
$ws.Range("K8").Value = 0.153344822805816
$ws.Range("L8").Value = 0.5042301998987817
$ws.Range("M8").Value = 0.6197304852958274
$ws.Range("N8").Value = 0.3258471028357235
$ws.Range("O8").Value = -0.005762179461639035
$ws.Range("P8").Value = 0.109569831889264
$ws.Range("Q8").Value = 0.1903397455737963
$ws.Range("R8").Value = 0.2565303361096417
$ws.Range("B10").Value = -0.05824571609997801
$ws.Range("C10").Value = -0.02294459231763052
$ws.Range("D10").Value = -0.3033481172167592
$ws.Range("E10").Value = 0.3816964293772955
$ws.Range("F10").Value = -0.4261300346409569
$ws.Range("G10").Value = 0.2529850742589501
$ws.Range("H10").Value = 0.9668095370549465
$ws.Range("K10").Value = 0.1435064096221038
$ws.Range("L10").Value = 0.5571744850310533
$ws.Range("M10").Value = 0.6057457838870046
$ws.Range("N10").Value = 0.431065527443969
$ws.Range("O10").Value = 0.02141133100128086
$ws.Range("P10").Value = 0.2280204233768813
$ws.Range("Q10").Value = 0.1929807552993139
$ws.Range("R10").Value = 0.3189536269809229
$ws.Range("B11").Value = -0.1494911127189555
$ws.Range("C11").Value = 0.158775156948585
$ws.Range("D11").Value = 0.4015523515806572
$ws.Range("E11").Value = -0.6093742791727917
$ws.Range("F11").Value = 0.3936451685116808
$ws.Range("G11").Value = 0.3350808734353892
$ws.Range("H11").Value = 0.153344822805816
$ws.Range("J11").Value = 0.1435064096221038
$ws.Range("L11").Value = -0.2630841920543963
$ws.Range("M11").Value = 0.2403457624510068
$ws.Range("N11").Value = -0.2172827789902214
$ws.Range("O11").Value = -0.2768602018790545
$ws.Range("P11").Value = -0.2261230290188869
$ws.Range("Q11").Value = 0.9212081916131816
$ws.Range("R11").Value = -0.362333011916985
$ws.Range("B12").Value = -0.1043516007832907
$ws.Range("C12").Value = 0.09866967570445913
$ws.Range("D12").Value = -0.5727500228386581
$ws.Range("E12").Value = 0.6801099074534445
$ws.Range("F12").Value = -0.694086968412218
$ws.Range("G12").Value = -0.08643860460969868
$ws.Range("H12").Value = 0.5042301998987817
$ws.Range("J12").Value = 0.5571744850310533
$ws.Range("K12").Value = -0.2630841920543963
$ws.Range("M12").Value = 0.2843010324439899
$ws.Range("N12").Value = 0.7871503025869309
$ws.Range("O12").Value = 0.07965621668814724
$ws.Range("P12").Value = 0.6670077063190093
$ws.Range("Q12").Value = -0.2304926056707239
$ws.Range("R12").Value = 0.8809551633059551
$ws.Range("B13").Value = -0.09656379635806157
$ws.Range("C13").Value = -0.04734339998996136
$ws.Range("D13").Value = -0.07574937763848384
$ws.Range("E13").Value = 0.07721240671536068
$ws.Range("F13").Value = -0.09685851345657408
$ws.Range("G13").Value = 0.779448656064079
$ws.Range("H13").Value = 0.6197304852958274
$ws.Range("J13").Value = 0.6057457838870046
$ws.Range("K13").Value = 0.2403457624510068
$ws.Range("L13").Value = 0.2843010324439899
$ws.Range("N13").Value = 0.3827780702837997
$ws.Range("O13").Value = 0.03840697036814979
$ws.Range("P13").Value = 0.2428911207911757
$ws.Range("Q13").Value = 0.2863353458644733
$ws.Range("R13").Value = 0.1709379276053136
$ws.Range("B14").Value = 0.04217011538748306
$ws.Range("C14").Value = -0.09583277978816673
$ws.Range("D14").Value = -0.4971150376806725
$ws.Range("E14").Value = 0.4942970076750814
$ws.Range("F14").Value = -0.5526442553679262
$ws.Range("G14").Value = 0.05749213078032143
$ws.Range("H14").Value = 0.3258471028357235
$ws.Range("J14").Value = 0.431065527443969
$ws.Range("K14").Value = -0.2172827789902214
$ws.Range("L14").Value = 0.7871503025869309
$ws.Range("M14").Value = 0.3827780702837997
$ws.Range("O14").Value = 0.213183587903477
$ws.Range("P14").Value = 0.9576445789701028
$ws.Range("Q14").Value = -0.1389052493132368
$ws.Range("R14").Value = 0.8156307420748278
$ws.Range("B15").Value = -0.1144752272881613
$ws.Range("C15").Value = 0.1405959713080313
$ws.Range("D15").Value = 0.1114722979239853
$ws.Range("E15").Value = -0.02813934050008031
$ws.Range("F15").Value = 0.05383151016975546
$ws.Range("G15").Value = 0.1739025695239392
$ws.Range("H15").Value = -0.005762179461639035
$ws.Range("J15").Value = 0.02141133100128086
$ws.Range("K15").Value = -0.2768602018790545
$ws.Range("L15").Value = 0.07965621668814724
$ws.Range("M15").Value = 0.03840697036814979
$ws.Range("N15").Value = 0.213183587903477
$ws.Range("P15").Value = 0.2407003385657958
$ws.Range("Q15").Value = -0.2164584615517044
$ws.Range("R15").Value = 0.02799107386774778
$ws.Range("B16").Value = 0.02034454368830499
$ws.Range("C16").Value = -0.04777815539842905
$ws.Range("D16").Value = -0.3819033532437513
$ws.Range("E16").Value = 0.3441744846429838
$ws.Range("F16").Value = -0.4195533350460475
$ws.Range("G16").Value = 0.02311440040280032
$ws.Range("H16").Value = 0.109569831889264
$ws.Range("J16").Value = 0.2280204233768813
$ws.Range("K16").Value = -0.2261230290188869
$ws.Range("L16").Value = 0.6670077063190093
$ws.Range("M16").Value = 0.2428911207911757
$ws.Range("N16").Value = 0.9576445789701028
$ws.Range("O16").Value = 0.2407003385657958
$ws.Range("Q16").Value = -0.1577812529357381
$ws.Range("R16").Value = 0.7654879483954553
$ws.Range("B17").Value = -0.143378829685207
$ws.Range("C17").Value = 0.1309388383387151
$ws.Range("D17").Value = 0.3454478151133055
$ws.Range("E17").Value = -0.550566853378591
$ws.Range("F17").Value = 0.3596911983350313
$ws.Range("G17").Value = 0.2711932158332709
$ws.Range("H17").Value = 0.1903397455737963
$ws.Range("J17").Value = 0.1929807552993139
$ws.Range("K17").Value = 0.9212081916131816
$ws.Range("L17").Value = -0.2304926056707239
$ws.Range("M17").Value = 0.2863353458644733
$ws.Range("N17").Value = -0.1389052493132368
$ws.Range("O17").Value = -0.2164584615517044
$ws.Range("P17").Value = -0.1577812529357381
$ws.Range("R17").Value = -0.3483315260529886
$ws.Range("B18").Value = 0.03120943056409502
$ws.Range("C18").Value = 0.01224759455026274
$ws.Range("D18").Value = -0.5891277964904179
$ws.Range("E18").Value = 0.6426805496216214
$ws.Range("F18").Value = -0.5916236419409803
$ws.Range("G18").Value = -0.1406704254457551
$ws.Range("H18").Value = 0.2565303361096417
$ws.Range("J18").Value = 0.3189536269809229
$ws.Range("K18").Value = -0.362333011916985
$ws.Range("L18").Value = 0.8809551633059551
$ws.Range("M18").Value = 0.1709379276053136
$ws.Range("N18").Value = 0.8156307420748278
$ws.Range("O18").Value = 0.02799107386774778
$ws.Range("P18").Value = 0.7654879483954553
$ws.Range("Q18").Value = -0.3483315260529886
